# Add 2022-Q1 data:
#  - rename the existing "总计" sheet (2nd sheet) to "2022-Q1" and replace its
#    content with the fund-holdings detail for 2022-Q1
#  - insert a brand-new "总计" sheet after it, rebuilding the summary table
#    with a 2022-Q1 row on top of the existing 2021-Q4 row

$wb = $excel.ActiveWorkbook

$detailSheet = $wb.Worksheets.Item(1)      # "2021-Q4" detail sheet (template for formatting)
$q1Sheet     = $wb.Worksheets.Item(2)      # currently named "总计"; becomes "2022-Q1"

# --- Step 1: turn the old "总计" sheet into the new "2022-Q1" detail sheet ---
$q1Sheet.Cells.Clear()
$q1Sheet.Name = "2022-Q1"

# Copy header row formatting (B1:H1) from the 2021-Q4 sheet so the new sheet
# keeps the same bold/centered/bordered look.
$detailSheet.Range("B1:H1").Copy($q1Sheet.Range("B1:H1"))

$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Cells that hold numeric-looking text must be pre-formatted as Text so the
# values round-trip as strings (matching how the source data was produced)
# instead of being auto-parsed into numbers.
$q1Sheet.Range("B2:B5").NumberFormat = "@"
$q1Sheet.Range("D2:G5").NumberFormat = "@"

$q1Data = @(
    @{ Idx = 0; Code = "320003"; Name = "诺安先锋混合";             Scale = "45.79"; Position = "69.96"; Ratio = "3.36"; Value = "1.5385"; Rank = 6 },
    @{ Idx = 1; Code = "001743"; Name = "诺安优选回报灵活配置混合"; Scale = "6.13";  Position = "71.32"; Ratio = "4.55"; Value = "0.2789"; Rank = 6 },
    @{ Idx = 2; Code = "160921"; Name = "大成多策略混合(LOF)";      Scale = "1.13";  Position = "79.19"; Ratio = "6.18"; Value = "0.0698"; Rank = 3 },
    @{ Idx = 3; Code = "001744"; Name = "诺安进取回报灵活配置混合"; Scale = "0.04";  Position = "62.10"; Ratio = "4.73"; Value = "0.0019"; Rank = 3 }
)

$row = 2
foreach ($item in $q1Data) {
    $q1Sheet.Range("A$row").Value = $item.Idx
    $detailSheet.Range("A2").Copy()
    $q1Sheet.Range("A$row").PasteSpecial(-4122)
    $q1Sheet.Range("B$row").Value = $item.Code
    $q1Sheet.Range("C$row").Value = $item.Name
    $q1Sheet.Range("D$row").Value = $item.Scale
    $q1Sheet.Range("E$row").Value = $item.Position
    $q1Sheet.Range("F$row").Value = $item.Ratio
    $q1Sheet.Range("G$row").Value = $item.Value
    $q1Sheet.Range("H$row").Value = $item.Rank
    $row++
}

# --- Step 2: insert a fresh "总计" sheet right after "2022-Q1" ---
# Duplicate the "2021-Q4" sheet (rather than Worksheets.Add) so the new sheet
# inherits the same sheetPr/pageMargins/outline settings used throughout the
# workbook, then wipe its contents and rebuild the summary table.
$detailSheet.Copy($null, $q1Sheet)
$totalSheet = $wb.Worksheets.Item(3)
$totalSheet.Cells.Clear()
$totalSheet.Name = "总计"

$detailSheet.Range("B1:D1").Copy($totalSheet.Range("B1:D1"))
$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$totalSheet.Range("A2").Value = 0
$detailSheet.Range("A2").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 1.89

$totalSheet.Range("A3").Value = 1
$detailSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 0.04

# Restore the original active sheet selection (2021-Q4, tab 0).
$detailSheet.Activate()
